$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Yes" from B2 (Status) to D2 (Data1)
$ws.Range("D2").Value = $ws.Range("B2").Value2
$ws.Range("B2").ClearContents()

# Remove "Yes" from B3 entirely (no replacement)
$ws.Range("B3").ClearContents()

# Update the selected cell to D2
$ws.Range("D2").Select()
